$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20

$data[0,0] = "ECs"
$data[0,1] = "B2m"
$data[0,2] = "Gm11127"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 1416.977457666667
$data[0,7] = 4250.932373
$data[0,8] = 0.1031800631271045
$data[0,9] = 0.1039170487194107
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.02669033333333333
$data[0,13] = 0.080071
$data[0,14] = 0.03151215526754497
$data[0,15] = 0.03167738918590879
$data[0,16] = 37.81960067094256
$data[0,17] = 340.376406038483
$data[0,18] = 0.003251426169776409
$data[0,19] = 0.003291820795335816
$data[1,0] = "ECs"
$data[1,1] = "B2m"
$data[1,2] = "Gm11127"
$data[1,3] = "Inflammatory-Mac"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 1416.977457666667
$data[1,7] = 4250.932373
$data[1,8] = 0.1031800631271045
$data[1,9] = 0.1039170487194107
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.5872193333333333
$data[1,13] = 1.761658
$data[1,14] = 0.6933051969416235
$data[1,15] = 0.6969405412505114
$data[1,16] = 832.0765580393816
$data[1,17] = 7488.689022354433
$data[1,18] = 0.07153527398678634
$data[1,19] = 0.07242400417966184
$data[2,0] = "ECs"
$data[2,1] = "B2m"
$data[2,2] = "Gm11127"
$data[2,3] = "MuSCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 1416.977457666667
$data[2,7] = 4250.932373
$data[2,8] = 0.1031800631271045
$data[2,9] = 0.1039170487194107
$data[2,10] = 1
$data[2,11] = 0.5
$data[2,12] = 0.013254
$data[2,13] = 0.026508
$data[2,14] = 0.01564844098048136
$data[2,15] = 0.01048699569806884
$data[2,16] = 18.780619223914
$data[2,17] = 112.683715343484
$data[2,18] = 0.001614607128206836
$data[2,19] = 0.00108977764287647
$data[3,0] = "ECs"
$data[3,1] = "B2m"
$data[3,2] = "Gm11127"
$data[3,3] = "Resolving-Mac"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 1416.977457666667
$data[3,7] = 4250.932373
$data[3,8] = 0.1031800631271045
$data[3,9] = 0.1039170487194107
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.2198216666666667
$data[3,13] = 0.659465
$data[3,14] = 0.2595342068103501
$data[3,15] = 0.2608950738655111
$data[3,16] = 311.4823463733828
$data[3,17] = 2803.341117360445
$data[3,18] = 0.02677875584233493
$data[3,19] = 0.02711144610153656
$data[4,0] = "FAPs"
$data[4,1] = "B2m"
$data[4,2] = "Gm11127"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1358.281941333333
$data[4,7] = 4074.845824
$data[4,8] = 0.09890603106838422
$data[4,9] = 0.09961248847575943
$data[4,10] = 1
$data[4,11] = 0.3333333333333333
$data[4,12] = 0.02669033333333333
$data[4,13] = 0.080071
$data[4,14] = 0.03151215526754497
$data[4,15] = 0.03167738918590879
$data[4,16] = 36.25299777483378
$data[4,17] = 326.276979973504
$data[4,18] = 0.00311674220792355
$data[4,19] = 0.003155463565223485
$data[5,0] = "FAPs"
$data[5,1] = "B2m"
$data[5,2] = "Gm11127"
$data[5,3] = "Inflammatory-Mac"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1358.281941333333
$data[5,7] = 4074.845824
$data[5,8] = 0.09890603106838422
$data[5,9] = 0.09961248847575943
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.5872193333333333
$data[5,13] = 1.761658
$data[5,14] = 0.6933051969416235
$data[5,15] = 0.6969405412505114
$data[5,16] = 797.6094160684657
$data[5,17] = 7178.484744616192
$data[5,18] = 0.06857206534858046
$data[5,19] = 0.06942398163360611
$data[6,0] = "FAPs"
$data[6,1] = "B2m"
$data[6,2] = "Gm11127"
$data[6,3] = "MuSCs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1358.281941333333
$data[6,7] = 4074.845824
$data[6,8] = 0.09890603106838422
$data[6,9] = 0.09961248847575943
$data[6,10] = 1
$data[6,11] = 0.5
$data[6,12] = 0.013254
$data[6,13] = 0.026508
$data[6,14] = 0.01564844098048136
$data[6,15] = 0.01048699569806884
$data[6,16] = 18.002668850432
$data[6,17] = 108.016013102592
$data[6,18] = 0.001547725189787266
$data[6,19] = 0.001044635738119221
$data[7,0] = "FAPs"
$data[7,1] = "B2m"
$data[7,2] = "Gm11127"
$data[7,3] = "Resolving-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1358.281941333333
$data[7,7] = 4074.845824
$data[7,8] = 0.09890603106838422
$data[7,9] = 0.09961248847575943
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.2198216666666667
$data[7,13] = 0.659465
$data[7,14] = 0.2595342068103501
$data[7,15] = 0.2608950738655111
$data[7,16] = 298.5798001471289
$data[7,17] = 2687.21820132416
$data[7,18] = 0.02566949832209294
$data[7,19] = 0.02598840753881063
$data[8,0] = "Inflammatory-Mac"
$data[8,1] = "B2m"
$data[8,2] = "Gm11127"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 5656.041829666666
$data[8,7] = 16968.125489
$data[8,8] = 0.4118560601490074
$data[8,9] = 0.4147978298401142
$data[8,10] = 1
$data[8,11] = 0.3333333333333333
$data[8,12] = 0.02669033333333333
$data[8,13] = 0.080071
$data[8,14] = 0.03151215526754497
$data[8,15] = 0.03167738918590879
$data[8,16] = 150.9616417810799
$data[8,17] = 1358.654776029719
$data[8,18] = 0.01297847211529486
$data[8,19] = 0.01313971228931567
$data[9,0] = "Inflammatory-Mac"
$data[9,1] = "B2m"
$data[9,2] = "Gm11127"
$data[9,3] = "Inflammatory-Mac"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 5656.041829666666
$data[9,7] = 16968.125489
$data[9,8] = 0.4118560601490074
$data[9,9] = 0.4147978298401142
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 0.5872193333333333
$data[9,13] = 1.761658
$data[9,14] = 0.6933051969416235
$data[9,15] = 0.6969405412505114
$data[9,16] = 3321.337112522306
$data[9,17] = 29892.03401270076
$data[9,18] = 0.2855419468932087
$data[9,19] = 0.2890894240383067
$data[10,0] = "Inflammatory-Mac"
$data[10,1] = "B2m"
$data[10,2] = "Gm11127"
$data[10,3] = "MuSCs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 5656.041829666666
$data[10,7] = 16968.125489
$data[10,8] = 0.4118560601490074
$data[10,9] = 0.4147978298401142
$data[10,10] = 1
$data[10,11] = 0.5
$data[10,12] = 0.013254
$data[10,13] = 0.026508
$data[10,14] = 0.01564844098048136
$data[10,15] = 0.01048699569806884
$data[10,16] = 74.96517841040199
$data[10,17] = 449.7910704624119
$data[10,18] = 0.006444905249695325
$data[10,19] = 0.004349983057101568
$data[11,0] = "Inflammatory-Mac"
$data[11,1] = "B2m"
$data[11,2] = "Gm11127"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 5656.041829666666
$data[11,7] = 16968.125489
$data[11,8] = 0.4118560601490074
$data[11,9] = 0.4147978298401142
$data[11,10] = 1
$data[11,11] = 0.3333333333333333
$data[11,12] = 0.2198216666666667
$data[11,13] = 0.659465
$data[11,14] = 0.2595342068103501
$data[11,15] = 0.2608950738655111
$data[11,16] = 1243.320541733709
$data[11,17] = 11189.88487560338
$data[11,18] = 0.1068907358908085
$data[11,19] = 0.1082187104553903
$data[12,0] = "MuSCs"
$data[12,1] = "B2m"
$data[12,2] = "Gm11127"
$data[12,3] = "ECs"
$data[12,4] = 2
$data[12,5] = 1
$data[12,6] = 292.1868055
$data[12,7] = 584.373611
$data[12,8] = 0.02127616983126987
$data[12,9] = 0.01428542627267643
$data[12,10] = 1
$data[12,11] = 0.3333333333333333
$data[12,12] = 0.02669033333333333
$data[12,13] = 0.080071
$data[12,14] = 0.03151215526754497
$data[12,15] = 0.03167738918590879
$data[12,16] = 7.798563234396833
$data[12,17] = 46.791379406381
$data[12,18] = 0.0006704579672216322
$data[12,19] = 0.0004525250077261775
$data[13,0] = "MuSCs"
$data[13,1] = "B2m"
$data[13,2] = "Gm11127"
$data[13,3] = "Inflammatory-Mac"
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 292.1868055
$data[13,7] = 584.373611
$data[13,8] = 0.02127616983126987
$data[13,9] = 0.01428542627267643
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 0.5872193333333333
$data[13,13] = 1.761658
$data[13,14] = 0.6933051969416235
$data[13,15] = 0.6969405412505114
$data[13,16] = 171.5777411345063
$data[13,17] = 1029.466446807038
$data[13,18] = 0.01475087911503198
$data[13,19] = 0.009956092718473383
$data[14,0] = "MuSCs"
$data[14,1] = "B2m"
$data[14,2] = "Gm11127"
$data[14,3] = "MuSCs"
$data[14,4] = 2
$data[14,5] = 1
$data[14,6] = 292.1868055
$data[14,7] = 584.373611
$data[14,8] = 0.02127616983126987
$data[14,9] = 0.01428542627267643
$data[14,10] = 1
$data[14,11] = 0.5
$data[14,12] = 0.013254
$data[14,13] = 0.026508
$data[14,14] = 0.01564844098048136
$data[14,15] = 0.01048699569806884
$data[14,16] = 3.872643920097
$data[14,17] = 15.490575680388
$data[14,18] = 0.0003329388878953246
$data[14,19] = 0.0001498112038666372
$data[15,0] = "MuSCs"
$data[15,1] = "B2m"
$data[15,2] = "Gm11127"
$data[15,3] = "Resolving-Mac"
$data[15,4] = 2
$data[15,5] = 1
$data[15,6] = 292.1868055
$data[15,7] = 584.373611
$data[15,8] = 0.02127616983126987
$data[15,9] = 0.01428542627267643
$data[15,10] = 1
$data[15,11] = 0.3333333333333333
$data[15,12] = 0.2198216666666667
$data[15,13] = 0.659465
$data[15,14] = 0.2595342068103501
$data[15,15] = 0.2608950738655111
$data[15,16] = 64.22899056301917
$data[15,17] = 385.373943378115
$data[15,18] = 0.005521893861120925
$data[15,19] = 0.003726997342610229
$data[16,0] = "Resolving-Mac"
$data[16,1] = "B2m"
$data[16,2] = "Gm11127"
$data[16,3] = "ECs"
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 5009.566731666667
$data[16,7] = 15028.700195
$data[16,8] = 0.3647816758242341
$data[16,9] = 0.3673872066920393
$data[16,10] = 1
$data[16,11] = 0.3333333333333333
$data[16,12] = 0.02669033333333333
$data[16,13] = 0.080071
$data[16,14] = 0.03151215526754497
$data[16,15] = 0.03167738918590879
$data[16,16] = 133.7070059237606
$data[16,17] = 1203.363053313845
$data[16,18] = 0.01149505680732852
$data[16,19] = 0.01163786752830764
$data[17,0] = "Resolving-Mac"
$data[17,1] = "B2m"
$data[17,2] = "Gm11127"
$data[17,3] = "Inflammatory-Mac"
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 5009.566731666667
$data[17,7] = 15028.700195
$data[17,8] = 0.3647816758242341
$data[17,9] = 0.3673872066920393
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 0.5872193333333333
$data[17,13] = 1.761658
$data[17,14] = 0.6933051969416235
$data[17,15] = 0.6969405412505114
$data[17,16] = 2941.714436458145
$data[17,17] = 26475.42992812331
$data[17,18] = 0.2529050315980161
$data[17,19] = 0.2560470386804634
$data[18,0] = "Resolving-Mac"
$data[18,1] = "B2m"
$data[18,2] = "Gm11127"
$data[18,3] = "MuSCs"
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = 5009.566731666667
$data[18,7] = 15028.700195
$data[18,8] = 0.3647816758242341
$data[18,9] = 0.3673872066920393
$data[18,10] = 1
$data[18,11] = 0.5
$data[18,12] = 0.013254
$data[18,13] = 0.026508
$data[18,14] = 0.01564844098048136
$data[18,15] = 0.01048699569806884
$data[18,16] = 66.39679746151
$data[18,17] = 398.38078476906
$data[18,18] = 0.005708264524896613
$data[18,19] = 0.003852788056104943
$data[19,0] = "Resolving-Mac"
$data[19,1] = "B2m"
$data[19,2] = "Gm11127"
$data[19,3] = "Resolving-Mac"
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 5009.566731666667
$data[19,7] = 15028.700195
$data[19,8] = 0.3647816758242341
$data[19,9] = 0.3673872066920393
$data[19,10] = 1
$data[19,11] = 0.3333333333333333
$data[19,12] = 0.2198216666666667
$data[19,13] = 0.659465
$data[19,14] = 0.2595342068103501
$data[19,15] = 0.2608950738655111
$data[19,16] = 1101.211308232853
$data[19,17] = 9910.901774095675
$data[19,18] = 0.09467332289399286
$data[19,19] = 0.09584951242716337

$ws.Range("A2:T21").Value = $data
